# Fixed a bug in chgSymbols
# Rewrites rows 2-19 and row 21 of the data table with the corrected
# row ordering/values (row 20 is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(101, 9, 30, 15, 60, 15)
    3  = @(1201, 2, 10, 10, 10, 10)
    4  = @(601, 9, 60, 67, 60, 42)
    5  = @(1203, 3, 15, 15, 15, 15)
    6  = @(902, 1, 0, 0, 0, 0)
    7  = @(701, 3, 90, 45, 97, 15)
    8  = @(401, 9, 48, 67, 75, 45)
    9  = @(801, 3, 67, 65, 52, 45)
    10 = @(301, 6, 45, 30, 60, 45)
    11 = @(501, 9, 52, 30, 75, 45)
    12 = @(201, 9, 30, 15, 45, 30)
    13 = @(1202, 2, 10, 10, 10, 10)
    14 = @(901, 16, 15, 45, 60, 60)
    15 = @(1001, 18, 30, 75, 60, 72)
    16 = @(3, 0, 3, 3, 3, 3)
    17 = @(1101, 0, 15, 30, 30, 0)
    18 = @(802, 0, 4, 5, 4, 0)
    19 = @(502, 0, 4, 0, 0, 0)
    21 = @(1, 0, 2, 2, 2, 2)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
